$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.14727665526569922
$ws.Range("B1").Value = 0.14719735498711373
$ws.Range("A2").Value = -0.14299875742206147
$ws.Range("B2").Value = 0.14267334027570922
$ws.Range("A3").Value = -0.092969991031120358
$ws.Range("B3").Value = 0.09278380273845066
$ws.Range("A4").Value = -0.084783802786516205
$ws.Range("B4").Value = 0.084219385823146808
$ws.Range("A5").Value = -0.081219385849396808
$ws.Range("B5").Value = 0.079287361564201575
$ws.Range("A6").Value = -0.011787657986337408
$ws.Range("B6").Value = 0.011454134889097745
$ws.Range("A7").Value = -0.0014541349564480832
$ws.Range("B7").Value = 0.0013788646266279514
$ws.Range("A8").Value = 0.0086211353059200135
$ws.Range("B8").Value = -0.0087456370230611924
$ws.Range("A9").Value = 0.010745636995562524
$ws.Range("B9").Value = -0.010844598989396825
$ws.Range("A10").Value = 0.012844598962812981
$ws.Range("B10").Value = -0.012850835331752464
$ws.Range("A11").Value = 0.01585083530030218
$ws.Range("B11").Value = -0.015861217118358972
$ws.Range("A12").Value = 0.018453317699624439
$ws.Range("B12").Value = -0.018537380147710447
$ws.Range("A13").Value = -0.0049335803672807543
$ws.Range("B13").Value = 0.0049236937944305836
$ws.Range("A14").Value = 0.0030763061521286161
$ws.Range("B14").Value = -0.0030764923911306141
$ws.Range("A15").Value = 0.0040764923721141599
$ws.Range("B15").Value = -0.0040774764540643105
$ws.Range("A16").Value = -0.0060337914763781697
$ws.Range("B16").Value = 0.0060034258686889608
$ws.Range("A17").Value = -0.004003425893126078
$ws.Range("B17").Value = 0.0039999999655497831
$ws.Range("A18").Value = 0.0018014024046948407
$ws.Range("B18").Value = -0.0018273818045813073
$ws.Range("A19").Value = 0.0058273817841012487
$ws.Range("B19").Value = -0.00600221109710386
$ws.Range("A20").Value = 0.010002211077022594
$ws.Range("B20").Value = -0.010043212489843967
$ws.Range("A21").Value = -0.0040056685496399069
$ws.Range("B21").Value = 0.0039999999799729125
$ws.Range("A22").Value = -0.045703349270414861
$ws.Range("B22").Value = 0.045492627277642583
$ws.Range("A23").Value = -0.040492627310005247
$ws.Range("B23").Value = 0.040097761112404484
$ws.Range("A24").Value = -0.02009776122144924
$ws.Range("B24").Value = 0.019999999889500408
$ws.Range("A25").Value = -0.065499703605706117
$ws.Range("B25").Value = 0.065414800926090422
$ws.Range("A26").Value = -0.062914800957326378
$ws.Range("B26").Value = 0.062807878159018671
$ws.Range("A27").Value = -0.060307878191612652
$ws.Range("B27").Value = 0.05969025777887893
$ws.Range("A28").Value = -0.0576902578146683
$ws.Range("B28").Value = 0.057281094548022438
$ws.Range("A29").Value = -0.050281094612701693
$ws.Range("B29").Value = 0.050172036716638502
$ws.Range("A30").Value = 0.0098279629565634075
$ws.Range("B30").Value = -0.0099435565435914697
$ws.Range("A31").Value = -0.014023139520130457
$ws.Range("B31").Value = 0.014001305860722013
$ws.Range("A32").Value = -0.0040013059417862706
$ws.Range("B32").Value = 0.0039999999482684956
